$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared-string-backed text cells
$ws.Range("E4").Value = "概率"
$ws.Range("F6").Value = "观看11-40次的档位（前后包含）"

# Update numeric values in the tiers table
$ws.Range("D5").Value = 10
$ws.Range("C6").Value = 11
$ws.Range("E7").Value = 0.003
$ws.Range("E8").Value = 0.005

# Update the active selection on the sheet
$ws.Range("F16").Select()

# Update the workbook window size
$excel.ActiveWindow.Width = 28800
$excel.ActiveWindow.Height = 12375
